# SKetched Manual.docx edits
# 1. "Team WaterBoard" + "Obv change this to real name" paragraphs -> single
#    paragraph "Team Root Tubers" (proofErr markers removed).
# 2. "<>" (credits placeholder) -> "Chase Allison".
# 3. The following empty paragraph gains a new "<>" run, and a brand new
#    empty paragraph is inserted right after it.
# 4. Two narrative paragraphs get their runs merged (proofErr markers
#    removed) while keeping the exact same visible text.

$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-RangeXml($range, [string]$bodyXml) {
    $full = $pkgHeader + '<w:body>' + $bodyXml + '</w:body>' + $pkgFooter
    $range.InsertXML($full)
}

# ---------------------------------------------------------------------
# Region 1: "Team WaterBoard" paragraph + "Obv change this to real name"
# paragraph collapse into a single paragraph "Team Root Tubers".
# ---------------------------------------------------------------------
$teamPara = $d.Paragraphs(11)
$obvPara = $d.Paragraphs(12)
$region1 = $d.Range($teamPara.Range.Start, $obvPara.Range.End)

$region1Xml = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Jumble" w:hAnsi="Jumble"/><w:sz w:val="60"/><w:szCs w:val="60"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Jumble" w:hAnsi="Jumble"/><w:sz w:val="60"/><w:szCs w:val="60"/></w:rPr><w:t xml:space="preserve">Team </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Jumble" w:hAnsi="Jumble"/><w:sz w:val="60"/><w:szCs w:val="60"/></w:rPr><w:t>Root Tubers</w:t></w:r>' + `
  '</w:p>'

Set-RangeXml $region1 $region1Xml

# ---------------------------------------------------------------------
# Region 2: "<>" -> "Chase Allison"; following empty paragraph gains a
# "<>" run; a brand new empty paragraph is inserted after it.
# ---------------------------------------------------------------------
$anglePara = $d.Paragraphs(13)
$blankPara = $d.Paragraphs(14)
$region2 = $d.Range($anglePara.Range.Start, $blankPara.Range.End)

$region2Xml = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Jumble" w:hAnsi="Jumble"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Jumble" w:hAnsi="Jumble"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>Chase Allison</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Jumble" w:hAnsi="Jumble"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Jumble" w:hAnsi="Jumble"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>&lt;&gt;</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Jumble" w:hAnsi="Jumble"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr></w:p>'

Set-RangeXml $region2 $region2Xml

# ---------------------------------------------------------------------
# Region 3: merge runs (drop proofErr markers) in the "An evil curse..."
# paragraph. Visible text is unchanged.
# ---------------------------------------------------------------------
$storyPara = $d.Paragraphs(21)
$region3 = $storyPara.Range

$region3Xml = '<w:p><w:pPr><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve">An evil curse has fallen upon </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>Lil’ Timmy</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>. You</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>, a water spirit,</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> have been summoned</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> by their Spirit Guardian, Karatian,</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> to protect </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve">Lil’ Timmy </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve">through any means possible. </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve">However, this curse is not the ordinary world-ending curses this Karatian usually deals with; this curse plagues Lil’ Timmy’s drawings. </w:t></w:r>' + `
  '</w:p>'

Set-RangeXml $region3 $region3Xml

# ---------------------------------------------------------------------
# Region 4: merge all runs (drop proofErr markers) in the "Karatian, with
# little experience..." paragraph. Visible text is unchanged.
# ---------------------------------------------------------------------
$goodluckPara = $d.Paragraphs(23)
$region4 = $goodluckPara.Range

$region4Xml = '<w:p><w:pPr><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>Karatian, with little experience with this curse, cast a faulty incantation, which left you in this infantile state, stripped of almost all of your normal powers. So, armed with only basic water balloons, you must fight your way through this curse and dispel it completely, all while keeping yourself alive. Good luck, and stay liquid!</w:t></w:r>' + `
  '</w:p>'

Set-RangeXml $region4 $region4Xml
